# Apply the "Updated cryptos list" refresh: new prices / volume deltas,
# plus the ImmutableX<->Algorand and FraxShare<->TrustWalletToken row swaps.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.165.47'
$ws.Range("D3").Value = '1.643.14'
$ws.Range("E3").Value = '  -2.52%  '
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = '  -0.41%  '
$ws.Range("D5").Value = "'307.13"
$ws.Range("E5").Value = '  -2.13%  '
$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = '  -0.31%  '
$ws.Range("D7").Value = "'0.3899"
$ws.Range("E7").Value = '  -0.90%  '
$ws.Range("D8").Value = "'0.3860"
$ws.Range("E8").Value = '  -2.93%  '
$ws.Range("D9").Value = "'0.9991"
$ws.Range("E9").Value = '  -0.58%  '
$ws.Range("D10").Value = "'49.64"
$ws.Range("E10").Value = '  -4.59%  '
$ws.Range("D11").Value = "'1.356"
$ws.Range("E11").Value = '  -5.10%  '
$ws.Range("D12").Value = "'0.08633"
$ws.Range("E12").Value = '  -0.62%  '
$ws.Range("D13").Value = "'23.65"
$ws.Range("E13").Value = '  -6.43%  '
$ws.Range("D14").Value = "'7.104"
$ws.Range("E14").Value = '  -3.10%  '
$ws.Range("E15").Value = '  -2.38%  '
$ws.Range("D16").Value = "'7.460"
$ws.Range("E16").Value = '  -4.46%  '
$ws.Range("D17").Value = '1.639.21'
$ws.Range("E17").Value = '  +0.13%  '
$ws.Range("D18").Value = "'94.76"
$ws.Range("E18").Value = '  +0.50%  '
$ws.Range("D19").Value = "'0.06898"
$ws.Range("E19").Value = '  -2.97%  '
$ws.Range("D20").Value = "'20.35"
$ws.Range("E20").Value = '  +0.96%  '
$ws.Range("D21").Value = "'6.895"
$ws.Range("E21").Value = '  -3.64%  '
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("D23").Value = "'13.58"
$ws.Range("E23").Value = '  -3.95%  '
$ws.Range("D24").Value = '24.167.67'
$ws.Range("D25").Value = "'2.391"
$ws.Range("E25").Value = '  +0.30%  '
$ws.Range("D26").Value = "'2.818"
$ws.Range("E26").Value = '  +1.12%  '
$ws.Range("E27").Value = '  -6.51%  '
$ws.Range("D28").Value = "'157.70"
$ws.Range("E28").Value = '  -3.05%  '
$ws.Range("D29").Value = "'8.568"
$ws.Range("E29").Value = '  +8.82%  '
$ws.Range("D30").Value = "'140.28"
$ws.Range("E30").Value = '  -7.21%  '
$ws.Range("D31").Value = "'5.355"
$ws.Range("E31").Value = '  -7.29%  '
$ws.Range("D32").Value = "'2.405"
$ws.Range("E32").Value = '  -7.81%  '
$ws.Range("D33").Value = '1.825.20'
$ws.Range("E33").Value = '  +0.33%  '
$ws.Range("D34").Value = "'7.008"
$ws.Range("E34").Value = '  +0.63%  '
$ws.Range("D35").Value = "'0.08083"
$ws.Range("D36").Value = "'0.02905"
$ws.Range("E36").Value = '  -5.96%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = "'0.9514"
$ws.Range("E37").Value = '  -6.68%  '
$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").Value = "'0.2679"
$ws.Range("E38").Value = '  -4.89%  '
$ws.Range("D39").Value = "'0.09200"
$ws.Range("E39").Value = '  -3.97%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = "'10.02"
$ws.Range("E40").Value = '  -4.81%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = "'1.459"
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").Value = "'0.7543"
$ws.Range("E42").Value = '  -5.59%  '
$ws.Range("D43").Value = "'13.01"
$ws.Range("E43").Value = '  -5.44%  '
$ws.Range("E44").Value = '  -4.27%  '
$ws.Range("D45").Value = "'0.6905"
$ws.Range("E45").Value = '  -4.14%  '
$ws.Range("E46").Value = '  -5.12%  '
$ws.Range("D47").Value = "'4.082"
$ws.Range("E47").Value = '  -2.44%  '
$ws.Range("D48").Value = "'0.9992"
$ws.Range("D49").Value = "'0.08402"
$ws.Range("E49").Value = '  -4.13%  '
$ws.Range("D50").Value = "'1.261"
$ws.Range("E50").Value = '  -6.24%  '
$ws.Range("D51").Value = "'133.24"
$ws.Range("E51").Value = '  -3.80%  '
